$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.6868273333333333
$ws.Cells.Item(2, 8).Value = 2.060482
$ws.Cells.Item(2, 9).Value = 0.01130642661970366
$ws.Cells.Item(2, 10).Value = 0.01130642661970366
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 63.55492266666666
$ws.Cells.Item(2, 14).Value = 190.664768
$ws.Cells.Item(2, 15).Value = 0.9936031556622397
$ws.Cells.Item(2, 16).Value = 0.9936031556622397
$ws.Cells.Item(2, 17).Value = 43.65125805535288
$ws.Cells.Item(2, 18).Value = 392.8613224981759
$ws.Cells.Item(2, 19).Value = 0.01123410116860111
$ws.Cells.Item(2, 20).Value = 0.01123410116860111

# Row 3
$ws.Cells.Item(3, 7).Value = 0.6868273333333333
$ws.Cells.Item(3, 8).Value = 2.060482
$ws.Cells.Item(3, 9).Value = 0.01130642661970366
$ws.Cells.Item(3, 10).Value = 0.01130642661970366
$ws.Cells.Item(3, 15).Value = 0.000612609346703606
$ws.Cells.Item(3, 16).Value = 0.000612609346703606
$ws.Cells.Item(3, 17).Value = 0.02691332905666667
$ws.Cells.Item(3, 18).Value = 0.24221996151
$ws.Cells.Item(3, 19).Value = 0.00000692642262504892
$ws.Cells.Item(3, 20).Value = 0.000006926422625048919

# Row 4
$ws.Cells.Item(4, 7).Value = 0.6868273333333333
$ws.Cells.Item(4, 8).Value = 2.060482
$ws.Cells.Item(4, 9).Value = 0.01130642661970366
$ws.Cells.Item(4, 10).Value = 0.01130642661970366
$ws.Cells.Item(4, 13).Value = 0.3699833333333333
$ws.Cells.Item(4, 14).Value = 1.10995
$ws.Cells.Item(4, 15).Value = 0.005784234991056675
$ws.Cells.Item(4, 16).Value = 0.005784234991056675
$ws.Cells.Item(4, 17).Value = 0.2541146662111111
$ws.Cells.Item(4, 18).Value = 2.2870319959
$ws.Cells.Item(4, 19).Value = 0.00006539902847750456
$ws.Cells.Item(4, 20).Value = 0.00006539902847750455

# Row 5
$ws.Cells.Item(5, 7).Value = 53.540432
$ws.Cells.Item(5, 9).Value = 0.8813728519762372
$ws.Cells.Item(5, 10).Value = 0.881372851976237
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 63.55492266666666
$ws.Cells.Item(5, 14).Value = 190.664768
$ws.Cells.Item(5, 15).Value = 0.9936031556622397
$ws.Cells.Item(5, 16).Value = 0.9936031556622397
$ws.Cells.Item(5, 17).Value = 3402.758015299925
$ws.Cells.Item(5, 18).Value = 30624.82213769932
$ws.Cells.Item(5, 19).Value = 0.8757348470386173
$ws.Cells.Item(5, 20).Value = 0.8757348470386171

# Row 6
$ws.Cells.Item(6, 7).Value = 53.540432
$ws.Cells.Item(6, 9).Value = 0.8813728519762372
$ws.Cells.Item(6, 10).Value = 0.881372851976237
$ws.Cells.Item(6, 15).Value = 0.000612609346703606
$ws.Cells.Item(6, 16).Value = 0.000612609346703606
$ws.Cells.Item(6, 19).Value = 0.0005399372470514566
$ws.Cells.Item(6, 20).Value = 0.0005399372470514566

# Row 7
$ws.Cells.Item(7, 7).Value = 53.540432
$ws.Cells.Item(7, 9).Value = 0.8813728519762372
$ws.Cells.Item(7, 10).Value = 0.881372851976237
$ws.Cells.Item(7, 13).Value = 0.3699833333333333
$ws.Cells.Item(7, 14).Value = 1.10995
$ws.Cells.Item(7, 15).Value = 0.005784234991056675
$ws.Cells.Item(7, 16).Value = 0.005784234991056675
$ws.Cells.Item(7, 17).Value = 19.80906749946667
$ws.Cells.Item(7, 18).Value = 178.2816074952
$ws.Cells.Item(7, 19).Value = 0.005098067690568367
$ws.Cells.Item(7, 20).Value = 0.005098067690568366

# Row 8
$ws.Cells.Item(8, 7).Value = 6.476716
$ws.Cells.Item(8, 8).Value = 19.430148
$ws.Cells.Item(8, 9).Value = 0.1066185206043934
$ws.Cells.Item(8, 10).Value = 0.1066185206043934
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 63.55492266666666
$ws.Cells.Item(8, 14).Value = 190.664768
$ws.Cells.Item(8, 15).Value = 0.9936031556622397
$ws.Cells.Item(8, 16).Value = 0.9936031556622397
$ws.Cells.Item(8, 17).Value = 411.6271845139626
$ws.Cells.Item(8, 18).Value = 3704.644660625664
$ws.Cells.Item(8, 19).Value = 0.1059364985245649
$ws.Cells.Item(8, 20).Value = 0.1059364985245649

# Row 9
$ws.Cells.Item(9, 7).Value = 6.476716
$ws.Cells.Item(9, 8).Value = 19.430148
$ws.Cells.Item(9, 9).Value = 0.1066185206043934
$ws.Cells.Item(9, 10).Value = 0.1066185206043934
$ws.Cells.Item(9, 15).Value = 0.000612609346703606
$ws.Cells.Item(9, 16).Value = 0.000612609346703606
$ws.Cells.Item(9, 17).Value = 0.25379011646
$ws.Cells.Item(9, 18).Value = 2.28411104814
$ws.Cells.Item(9, 19).Value = 0.00006531550225396243
$ws.Cells.Item(9, 20).Value = 0.00006531550225396243

# Row 10
$ws.Cells.Item(10, 7).Value = 6.476716
$ws.Cells.Item(10, 8).Value = 19.430148
$ws.Cells.Item(10, 9).Value = 0.1066185206043934
$ws.Cells.Item(10, 10).Value = 0.1066185206043934
$ws.Cells.Item(10, 13).Value = 0.3699833333333333
$ws.Cells.Item(10, 14).Value = 1.10995
$ws.Cells.Item(10, 15).Value = 0.005784234991056675
$ws.Cells.Item(10, 16).Value = 0.005784234991056675
$ws.Cells.Item(10, 17).Value = 2.396276974733333
$ws.Cells.Item(10, 18).Value = 21.5664927726
$ws.Cells.Item(10, 19).Value = 0.0006167065775746297
$ws.Cells.Item(10, 20).Value = 0.0006167065775746297

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.04265633333333333
$ws.Cells.Item(11, 8).Value = 0.127969
$ws.Cells.Item(11, 9).Value = 0.0007022007996657373
$ws.Cells.Item(11, 10).Value = 0.0007022007996657372
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 63.55492266666666
$ws.Cells.Item(11, 14).Value = 190.664768
$ws.Cells.Item(11, 15).Value = 0.9936031556622397
$ws.Cells.Item(11, 16).Value = 0.9936031556622397
$ws.Cells.Item(11, 17).Value = 2.711019966243555
$ws.Cells.Item(11, 18).Value = 24.399179696192
$ws.Cells.Item(11, 19).Value = 0.0006977089304564248
$ws.Cells.Item(11, 20).Value = 0.0006977089304564247

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.04265633333333333
$ws.Cells.Item(12, 8).Value = 0.127969
$ws.Cells.Item(12, 9).Value = 0.0007022007996657373
$ws.Cells.Item(12, 10).Value = 0.0007022007996657372
$ws.Cells.Item(12, 15).Value = 0.000612609346703606
$ws.Cells.Item(12, 16).Value = 0.000612609346703606
$ws.Cells.Item(12, 17).Value = 0.001671488421666667
$ws.Cells.Item(12, 18).Value = 0.015043395795
$ws.Cells.Item(12, 19).Value = 0.0000004301747731379771
$ws.Cells.Item(12, 20).Value = 0.000000430174773137977

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.04265633333333333
$ws.Cells.Item(13, 8).Value = 0.127969
$ws.Cells.Item(13, 9).Value = 0.0007022007996657373
$ws.Cells.Item(13, 10).Value = 0.0007022007996657372
$ws.Cells.Item(13, 13).Value = 0.3699833333333333
$ws.Cells.Item(13, 14).Value = 1.10995
$ws.Cells.Item(13, 15).Value = 0.005784234991056675
$ws.Cells.Item(13, 16).Value = 0.005784234991056675
$ws.Cells.Item(13, 17).Value = 0.01578213239444444
$ws.Cells.Item(13, 18).Value = 0.14203919155
$ws.Cells.Item(13, 19).Value = 0.000004061694436174536
$ws.Cells.Item(13, 20).Value = 0.000004061694436174536
